$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '29.961.35'
Set-TextValue $ws.Range('E2') '  +1.62%  '
Set-TextValue $ws.Range('D3') '1.942.16'
Set-TextValue $ws.Range('E3') '  +1.19%  '
Set-TextValue $ws.Range('D4') '1.011'
Set-TextValue $ws.Range('E4') '  +0.30%  '
Set-TextValue $ws.Range('D5') '335.23'
Set-TextValue $ws.Range('E5') '  +2.97%  '
Set-TextValue $ws.Range('D6') '1.010'
Set-TextValue $ws.Range('E6') '  +0.22%  '
Set-TextValue $ws.Range('D7') '0.4847'
Set-TextValue $ws.Range('E7') '  +0.40%  '
Set-TextValue $ws.Range('D8') '0.4153'
Set-TextValue $ws.Range('E8') '  +1.90%  '
Set-TextValue $ws.Range('D9') '0.08212'
Set-TextValue $ws.Range('E9') '  +0.35%  '
Set-TextValue $ws.Range('D10') '1.018'
Set-TextValue $ws.Range('E10') '  -0.29%  '
Set-TextValue $ws.Range('D11') '23.88'
Set-TextValue $ws.Range('E11') '  +1.71%  '
Set-TextValue $ws.Range('D12') '1.972.33'
Set-TextValue $ws.Range('E12') '  +2.63%  '
Set-TextValue $ws.Range('D13') '6.095'
Set-TextValue $ws.Range('E13') '  +0.90%  '
Set-TextValue $ws.Range('D14') '7.330'
Set-TextValue $ws.Range('E14') '  +1.44%  '
Set-TextValue $ws.Range('D15') '91.55'
Set-TextValue $ws.Range('E15') '  +0.33%  '
Set-TextValue $ws.Range('E16') '  +1.20%  '
Set-TextValue $ws.Range('D17') '1.012'
Set-TextValue $ws.Range('E17') '  +0.39%  '
Set-TextValue $ws.Range('D18') '0.00001040'
Set-TextValue $ws.Range('E18') '  +0.09%  '
Set-TextValue $ws.Range('D19') '17.89'
Set-TextValue $ws.Range('E19') '  +0.64%  '
Set-TextValue $ws.Range('D20') '1.010'
Set-TextValue $ws.Range('E20') '  +0.23%  '
Set-TextValue $ws.Range('D21') '29.951.32'
Set-TextValue $ws.Range('E21') '  +1.47%  '
Set-TextValue $ws.Range('D22') '5.651'
Set-TextValue $ws.Range('E22') '  +0.21%  '
Set-TextValue $ws.Range('D23') '11.93'
Set-TextValue $ws.Range('E23') '  +1.37%  '
Set-TextValue $ws.Range('D24') '2.190'
Set-TextValue $ws.Range('E24') '  -0.32%  '
Set-TextValue $ws.Range('D25') '2.189.05'
Set-TextValue $ws.Range('E25') '  +2.03%  '
Set-TextValue $ws.Range('D26') '6.584'
Set-TextValue $ws.Range('E26') '  -1.39%  '
Set-TextValue $ws.Range('D27') '157.11'
Set-TextValue $ws.Range('E27') '  +0.26%  '
Set-TextValue $ws.Range('E28') '  +0.19%  '
Set-TextValue $ws.Range('D29') '2.110'
Set-TextValue $ws.Range('E29') '  -0.13%  '
Set-TextValue $ws.Range('D30') '121.33'
Set-TextValue $ws.Range('E30') '  +0.92%  '
Set-TextValue $ws.Range('D31') '1.018'
Set-TextValue $ws.Range('E31') '  -0.37%  '
Set-TextValue $ws.Range('D32') '0.09633'
Set-TextValue $ws.Range('E32') '  +0.80%  '
Set-TextValue $ws.Range('D33') '5.633'
Set-TextValue $ws.Range('E33') '  +2.00%  '
Set-TextValue $ws.Range('D34') '1.421'
Set-TextValue $ws.Range('E34') '  +2.65%  '
Set-TextValue $ws.Range('D35') '3.565'
Set-TextValue $ws.Range('E35') '  +0.06%  '
Set-TextValue $ws.Range('D36') '0.06551'
Set-TextValue $ws.Range('E36') '  +6.83%  '
Set-TextValue $ws.Range('D37') '0.02294'
Set-TextValue $ws.Range('E37') '  +0.42%  '
Set-TextValue $ws.Range('D38') '1.220'
Set-TextValue $ws.Range('E38') '  +3.14%  '
Set-TextValue $ws.Range('D39') '0.5975'
Set-TextValue $ws.Range('E39') '  -0.10%  '
Set-TextValue $ws.Range('D40') '8.009'
Set-TextValue $ws.Range('E40') '  -0.36%  '
Set-TextValue $ws.Range('E41') '  -0.62%  '
Set-TextValue $ws.Range('D42') '2.543'
Set-TextValue $ws.Range('E42') '  +5.98%  '
Set-TextValue $ws.Range('D43') '0.1855'
Set-TextValue $ws.Range('E43') '  -0.08%  '
Set-TextValue $ws.Range('D44') '1.244'
Set-TextValue $ws.Range('E44') '  -3.05%  '
Set-TextValue $ws.Range('D45') '12.39'
Set-TextValue $ws.Range('E45') '  -0.08%  '
Set-TextValue $ws.Range('D46') '0.07523'
Set-TextValue $ws.Range('E46') '  -1.12%  '
Set-TextValue $ws.Range('D47') '0.5583'
Set-TextValue $ws.Range('E47') '  +0.16%  '
Set-TextValue $ws.Range('E48') '  +1.58%  '
Set-TextValue $ws.Range('D49') '117.70'
Set-TextValue $ws.Range('E49') '  +0.23%  '
Set-TextValue $ws.Range('D50') '2.435'
Set-TextValue $ws.Range('E50') '  +0.11%  '
Set-TextValue $ws.Range('D51') '72.96'
Set-TextValue $ws.Range('E51') '  +0.40%  '
